$wb = $excel.ActiveWorkbook

# ALC row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 397.25
$ws.Range("I12").Value = 397.25
$ws.Range("K12").Value = 397.25
$ws.Range("M12").Value = -227.25

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3256
$ws.Range("J70").Value = 3292.5715
$ws.Range("L70").Value = 9877.7145
$ws.Range("N70").Value = -10417.7145

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 3256
$ws.Range("J73").Value = 3292.5715
$ws.Range("L73").Value = 9877.7145
$ws.Range("N73").Value = -11749.7145

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 916
$ws.Range("I76").Value = 949
$ws.Range("J76").Value = 850
$ws.Range("K76").Value = 949
$ws.Range("L76").Value = 850
$ws.Range("M76").Value = -634
$ws.Range("N76").Value = -1480

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 916
$ws.Range("I79").Value = 949
$ws.Range("J79").Value = 850
$ws.Range("K79").Value = 949
$ws.Range("L79").Value = 850
$ws.Range("M79").Value = 143
$ws.Range("N79").Value = -3034

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4954.778
$ws.Range("I116").Value = 3866
$ws.Range("J116").Value = 5499.1665
$ws.Range("K116").Value = 3866
$ws.Range("L116").Value = 5499.1665
$ws.Range("M116").Value = -424
$ws.Range("N116").Value = -12383.1665

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 9682.375
$ws.Range("I125").Value = 1449.5
$ws.Range("K125").Value = 13045.5
$ws.Range("M125").Value = -10585.5

# ARM row 46
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 4910.4
$ws.Range("J46").Value = 4910.4
$ws.Range("L46").Value = 4910.4
$ws.Range("N46").Value = -5548.4

# ARM row 57
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H57").Value = 15000
$ws.Range("I57").Value = 15000
$ws.Range("K57").Value = 15000
$ws.Range("M57").Value = -14516

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4000
$ws.Range("I61").Value = 4000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3788
$ws.Range("N61").ClearContents()

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2947
$ws.Range("I74").Value = 2914.8
$ws.Range("K74").Value = 2914.8
$ws.Range("M74").Value = -2040.8

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2947
$ws.Range("I77").Value = 2914.8
$ws.Range("K77").Value = 14574
$ws.Range("M77").Value = -10206

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 4625985
$ws.Range("I110").Value = 7400799
$ws.Range("K110").Value = 7400799
$ws.Range("M110").Value = -7398754

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1316.5714
$ws.Range("I122").Value = 1119.3334
$ws.Range("K122").Value = 3358.0002
$ws.Range("M122").Value = -908.0001999999999

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4000
$ws.Range("I136").Value = 4000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 12000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -9450
$ws.Range("N136").ClearContents()

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3055
$ws.Range("I99").Value = 1199
$ws.Range("K99").Value = 1199
$ws.Range("M99").Value = 299

# BSM row 126
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# BSM row 128
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 8000.5
$ws.Range("I128").Value = 8000.5
$ws.Range("K128").Value = 24001.5
$ws.Range("M128").Value = -21511.5

# CRP row 6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 200290.4
$ws.Range("I6").Value = 250275.25
$ws.Range("J6").Value = 351
$ws.Range("K6").Value = 250275.25
$ws.Range("L6").Value = 351
$ws.Range("M6").Value = -250162.25
$ws.Range("N6").Value = -577

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2531.5715
$ws.Range("I31").Value = 1945.9
$ws.Range("J31").Value = 3995.75
$ws.Range("K31").Value = 1945.9
$ws.Range("L31").Value = 3995.75
$ws.Range("M31").Value = -1650.9
$ws.Range("N31").Value = -4585.75

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2531.5715
$ws.Range("I34").Value = 1945.9
$ws.Range("J34").Value = 3995.75
$ws.Range("K34").Value = 1945.9
$ws.Range("L34").Value = 3995.75
$ws.Range("M34").Value = -1743.9
$ws.Range("N34").Value = -4399.75

# CUL row 33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2000
$ws.Range("J39").Value = 2000
$ws.Range("L39").Value = 6000
$ws.Range("N39").Value = -6588

# CUL row 44
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 425
$ws.Range("I44").Value = 173.75
$ws.Range("K44").Value = 521.25
$ws.Range("M44").Value = -123.25

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 5661.3335
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

# CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 3998
$ws.Range("I98").Value = 3994.5
$ws.Range("K98").Value = 11983.5
$ws.Range("M98").Value = -10485.5

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1433.7142
$ws.Range("J113").Value = 1523
$ws.Range("L113").Value = 4569
$ws.Range("N113").Value = -8909

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 16565
$ws.Range("I121").Value = 27596.25
$ws.Range("J121").Value = 5533.75
$ws.Range("K121").Value = 82788.75
$ws.Range("L121").Value = 16601.25
$ws.Range("M121").Value = -81478.75
$ws.Range("N121").Value = -19221.25

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3055.3333
$ws.Range("J137").Value = 3694.182
$ws.Range("L137").Value = 11082.546
$ws.Range("N137").Value = -21282.546

# CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2737
$ws.Range("I139").Value = 2737
$ws.Range("K139").Value = 8211
$ws.Range("M139").Value = -3071

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1532.5
$ws.Range("I80").Value = 1698.75
$ws.Range("K80").Value = 1698.75
$ws.Range("M80").Value = -700.75

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1532.5
$ws.Range("I83").Value = 1698.75
$ws.Range("K83").Value = 8493.75
$ws.Range("M83").Value = -3501.75

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 999.1

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4441.4546
$ws.Range("I122").Value = 4675.6
$ws.Range("J122").Value = 2100
$ws.Range("K122").Value = 14026.8
$ws.Range("L122").Value = 6300
$ws.Range("M122").Value = -11576.8
$ws.Range("N122").Value = -11200

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 8623.450000000001
$ws.Range("I132").Value = 8964.6
$ws.Range("K132").Value = 26893.8
$ws.Range("M132").Value = -24363.8

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9505.875
$ws.Range("I7").Value = 12668.333
$ws.Range("K7").Value = 12668.333
$ws.Range("M7").Value = -12556.333

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4803.5835
$ws.Range("I40").Value = 4614.3
$ws.Range("K40").Value = 4614.3
$ws.Range("M40").Value = -4478.3

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1619
$ws.Range("J46").Value = 1269.2
$ws.Range("L46").Value = 1269.2
$ws.Range("N46").Value = -1645.2

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1900
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# LTW row 107
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 12500
$ws.Range("I107").Value = 12500
$ws.Range("K107").Value = 12500
$ws.Range("M107").Value = -10580

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1900
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 8785.700000000001
$ws.Range("I122").Value = 9553.25
$ws.Range("J122").Value = 8274
$ws.Range("K122").Value = 28659.75
$ws.Range("L122").Value = 24822
$ws.Range("M122").Value = -26209.75
$ws.Range("N122").Value = -29722

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 9505.875
$ws.Range("I126").Value = 12668.333
$ws.Range("K126").Value = 38004.999
$ws.Range("M126").Value = -35534.999

# WVR row 46
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 771203.9399999999
$ws.Range("I81").Value = 1593.1428
$ws.Range("J81").Value = 1669083.1
$ws.Range("K81").Value = 3186.2856
$ws.Range("L81").Value = 3338166.2
$ws.Range("M81").Value = -2125.2856
$ws.Range("N81").Value = -3340288.2

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 771203.9399999999
$ws.Range("I84").Value = 1593.1428
$ws.Range("J84").Value = 1669083.1
$ws.Range("K84").Value = 15931.428
$ws.Range("L84").Value = 16690831
$ws.Range("M84").Value = -10627.428
$ws.Range("N84").Value = -16701439

# WVR row 97
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 572
$ws.Range("J97").Value = 572
$ws.Range("L97").Value = 572
$ws.Range("N97").Value = -2554

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2937.5
$ws.Range("I107").Value = 3016.6667
$ws.Range("K107").Value = 9050.000100000001
$ws.Range("M107").Value = -7130.000100000001

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2470.4285
$ws.Range("I126").Value = 998.5
$ws.Range("J126").Value = 4433
$ws.Range("K126").Value = 2995.5
$ws.Range("L126").Value = 13299
$ws.Range("M126").Value = -525.5
$ws.Range("N126").Value = -18239

# WVR row 134
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
